# Auto-generated edit script applying the Cerberus_Profits diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: H3 0->20607, J3 0->20607, L3 0->20607, N3 None->-20835
$ws.Range("H3").Value = 20607
$ws.Range("J3").Value = 20607
$ws.Range("L3").Value = 20607
$ws.Range("N3").Value = -20835

# Row 32: H32 2967->2580.2, J32 2967->2580.2, L32 2967->2580.2, N32 -3619->-3232.2
$ws.Range("H32").Value = 2580.2
$ws.Range("J32").Value = 2580.2
$ws.Range("L32").Value = 2580.2
$ws.Range("N32").Value = -3232.2

# Row 102: H102 0->20607, J102 0->20607, L102 0->20607, N102 None->-27097
$ws.Range("H102").Value = 20607
$ws.Range("J102").Value = 20607
$ws.Range("L102").Value = 20607
$ws.Range("N102").Value = -27097

# Row 112: H112 3256.861->3232.8572, J112 4163.654->4166.32, L112 12490.962->12498.96, N112 -14706.962->-14714.96
$ws.Range("H112").Value = 3232.8572
$ws.Range("J112").Value = 4166.32
$ws.Range("L112").Value = 12498.96
$ws.Range("N112").Value = -14714.96

# Row 137: H137 1648.5106->1615.3959, I137 1518.7028->1480.2894, K137 4556.1084->4440.8682, M137 -2006.1084->-1890.8682
$ws.Range("H137").Value = 1615.3959
$ws.Range("I137").Value = 1480.2894
$ws.Range("K137").Value = 4440.8682
$ws.Range("M137").Value = -1890.8682

# Row 138: H138 3022.476->2720.5, J138 2853.7297->2426.8542, L138 8561.1891->7280.562600000001, N138 -18841.1891->-17560.5626
$ws.Range("H138").Value = 2720.5
$ws.Range("J138").Value = 2426.8542
$ws.Range("L138").Value = 7280.562600000001
$ws.Range("N138").Value = -17560.5626

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32 3049.3823->3881.6155, I32 2215.3438->2838.9167, K32 2215.3438->2838.9167, M32 -1928.3438->-2551.9167
$ws.Range("H32").Value = 3881.6155
$ws.Range("I32").Value = 2838.9167
$ws.Range("K32").Value = 2838.9167
$ws.Range("M32").Value = -2551.9167

# Row 74: H74 1380.3636->1280.9323, I74 646.1515000000001->566.9729599999999, K74 646.1515000000001->566.9729599999999, M74 227.8484999999999->307.0270400000001
$ws.Range("H74").Value = 1280.9323
$ws.Range("I74").Value = 566.9729599999999
$ws.Range("K74").Value = 566.9729599999999
$ws.Range("M74").Value = 307.0270400000001

# Row 77: H77 1380.3636->1280.9323, I77 646.1515000000001->566.9729599999999, K77 3230.7575->2834.8648, M77 1137.2425->1533.1352
$ws.Range("H77").Value = 1280.9323
$ws.Range("I77").Value = 566.9729599999999
$ws.Range("K77").Value = 2834.8648
$ws.Range("M77").Value = 1533.1352

$ws = $wb.Worksheets.Item("BSM")
# Row 20: H20 1162.3077->1064.3334, I20 1219.091->1068.9286, J20 850->1000, K20 1219.091->1068.9286, L20 850->1000, M20 -972.0909999999999->-821.9286, N20 -1344->-1494
$ws.Range("H20").Value = 1064.3334
$ws.Range("I20").Value = 1068.9286
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1068.9286
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -821.9286
$ws.Range("N20").Value = -1494

# Row 94: H94 5991.5->8479.075000000001, I94 4404.615->4560.92, J94 8938.571->15009.333, K94 4404.615->4560.92, L94 8938.571->15009.333, M94 -3953.615->-4109.92, N94 -9840.571->-15911.333
$ws.Range("H94").Value = 8479.075000000001
$ws.Range("I94").Value = 4560.92
$ws.Range("J94").Value = 15009.333
$ws.Range("K94").Value = 4560.92
$ws.Range("L94").Value = 15009.333
$ws.Range("M94").Value = -4109.92
$ws.Range("N94").Value = -15911.333

# Row 103: H103 31462.5->30666, J103 31462.5->30666, L103 31462.5->30666, N103 -33806.5->-33010
$ws.Range("H103").Value = 30666
$ws.Range("J103").Value = 30666
$ws.Range("L103").Value = 30666
$ws.Range("N103").Value = -33010

$ws = $wb.Worksheets.Item("CRP")
# Row 33: H33 0->300, I33 0->300, K33 0->300, M33 None->79
$ws.Range("H33").Value = 300
$ws.Range("I33").Value = 300
$ws.Range("K33").Value = 300
$ws.Range("M33").Value = 79

# Row 86: H86 4723.6113->4525.2856, I86 4230.5557->4043.818, J86 5216.6665->5054.9, K86 4230.5557->4043.818, L86 5216.6665->5054.9, M86 -3107.5557->-2920.818, N86 -7462.6665->-7300.9
$ws.Range("H86").Value = 4525.2856
$ws.Range("I86").Value = 4043.818
$ws.Range("J86").Value = 5054.9
$ws.Range("K86").Value = 4043.818
$ws.Range("L86").Value = 5054.9
$ws.Range("M86").Value = -2920.818
$ws.Range("N86").Value = -7300.9

# Row 89: H89 4723.6113->4525.2856, I89 4230.5557->4043.818, J89 5216.6665->5054.9, K89 21152.7785->20219.09, L89 26083.3325->25274.5, M89 -15536.7785->-14603.09, N89 -37315.3325->-36506.5
$ws.Range("H89").Value = 4525.2856
$ws.Range("I89").Value = 4043.818
$ws.Range("J89").Value = 5054.9
$ws.Range("K89").Value = 20219.09
$ws.Range("L89").Value = 25274.5
$ws.Range("M89").Value = -14603.09
$ws.Range("N89").Value = -36506.5

# Row 116: H116 90000->0, J116 90000->0, L116 90000->0, N116 -99178->None
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 132: H132 4036.484->3782.6177, I132 1963.7826->1870.9615, K132 5891.3478->5612.8845, M132 -3361.3478->-3082.8845
$ws.Range("H132").Value = 3782.6177
$ws.Range("I132").Value = 1870.9615
$ws.Range("K132").Value = 5612.8845
$ws.Range("M132").Value = -3082.8845

# Row 134: H134 1161.3125->1121.22, I134 1046.2046->1007.63043, K134 3138.6138->3022.89129, M134 -603.6138000000001->-487.89129
$ws.Range("H134").Value = 1121.22
$ws.Range("I134").Value = 1007.63043
$ws.Range("K134").Value = 3022.89129
$ws.Range("M134").Value = -487.89129

# Row 141: H141 225796.08->201855.92, J141 259255.4->225748.67, L141 259255.4->225748.67, N141 -269615.4->-236108.67
$ws.Range("H141").Value = 201855.92
$ws.Range("J141").Value = 225748.67
$ws.Range("L141").Value = 225748.67
$ws.Range("N141").Value = -236108.67

$ws = $wb.Worksheets.Item("CUL")
# Row 36: H36 98->0, I36 98->0, K36 294->0, M36 -125->None
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Row 39: H39 6319.6->6299.6, J39 15499->10366, L39 46497->31098, N39 -47085->-31686
$ws.Range("H39").Value = 6299.6
$ws.Range("J39").Value = 10366
$ws.Range("L39").Value = 31098
$ws.Range("N39").Value = -31686

# Row 63: H63 2500->3000, I63 2500->3000, K63 7500->9000, M63 -6751->-8251
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 9000
$ws.Range("M63").Value = -8251

# Row 64: H64 3490347.5->2328565, I64 0->5000, K64 0->15000, M64 None->-14730
$ws.Range("H64").Value = 2328565
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 15000
$ws.Range("M64").Value = -14730

# Row 66: H66 2500->3000, I66 2500->3000, K66 22500->27000, M66 -18756->-23256
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 27000
$ws.Range("M66").Value = -23256

# Row 67: H67 3490347.5->2328565, I67 0->5000, K67 0->15000, M67 None->-14064
$ws.Range("H67").Value = 2328565
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 15000
$ws.Range("M67").Value = -14064

# Row 70: H70 1943.9166->1647.8422, I70 1200.5->1078.25, J70 2092.6->1799.7333, K70 3601.5->3234.75, L70 6277.799999999999->5399.199900000001, M70 -3286.5->-2919.75, N70 -6907.799999999999->-6029.199900000001
$ws.Range("H70").Value = 1647.8422
$ws.Range("I70").Value = 1078.25
$ws.Range("J70").Value = 1799.7333
$ws.Range("K70").Value = 3234.75
$ws.Range("L70").Value = 5399.199900000001
$ws.Range("M70").Value = -2919.75
$ws.Range("N70").Value = -6029.199900000001

# Row 73: H73 1943.9166->1647.8422, I73 1200.5->1078.25, J73 2092.6->1799.7333, K73 3601.5->3234.75, L73 6277.799999999999->5399.199900000001, M73 -2509.5->-2142.75, N73 -8461.799999999999->-7583.199900000001
$ws.Range("H73").Value = 1647.8422
$ws.Range("I73").Value = 1078.25
$ws.Range("J73").Value = 1799.7333
$ws.Range("K73").Value = 3234.75
$ws.Range("L73").Value = 5399.199900000001
$ws.Range("M73").Value = -2142.75
$ws.Range("N73").Value = -7583.199900000001

# Row 92: H92 347.65384->347.26923, I92 79.22727->78.77273, K92 237.68181->236.31819, M92 1010.31819->1011.68181
$ws.Range("H92").Value = 347.26923
$ws.Range("I92").Value = 78.77273
$ws.Range("K92").Value = 236.31819
$ws.Range("M92").Value = 1011.68181

# Row 130: H130 26999.875->23499.857, I130 0->8500, J130 26999.875->25999.834, K130 0->25500, L130 80999.625->77999.50199999999, M130 None->-20480, N130 -91039.625->-88039.50199999999
$ws.Range("H130").Value = 23499.857
$ws.Range("I130").Value = 8500
$ws.Range("J130").Value = 25999.834
$ws.Range("K130").Value = 25500
$ws.Range("L130").Value = 77999.50199999999
$ws.Range("M130").Value = -20480
$ws.Range("N130").Value = -88039.50199999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70 12056.363->11310.077, J70 8157->7885.143, L70 8157->7885.143, N70 -8697->-8425.143
$ws.Range("H70").Value = 11310.077
$ws.Range("J70").Value = 7885.143
$ws.Range("L70").Value = 7885.143
$ws.Range("N70").Value = -8425.143

# Row 73: H73 12056.363->11310.077, J73 8157->7885.143, L73 8157->7885.143, N73 -10029->-9757.143
$ws.Range("H73").Value = 11310.077
$ws.Range("J73").Value = 7885.143
$ws.Range("L73").Value = 7885.143
$ws.Range("N73").Value = -9757.143

# Row 102: H102 11765.173->11396.333, J102 4068.8->3858.25, L102 4068.8->3858.25, N102 -7312.8->-7102.25
$ws.Range("H102").Value = 11396.333
$ws.Range("J102").Value = 3858.25
$ws.Range("L102").Value = 3858.25
$ws.Range("N102").Value = -7102.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16 617.1429000000001->701.7778, I16 345.35294->395.92856, K16 345.35294->395.92856, M16 -175.35294->-225.92856
$ws.Range("H16").Value = 701.7778
$ws.Range("I16").Value = 395.92856
$ws.Range("K16").Value = 395.92856
$ws.Range("M16").Value = -225.92856

# Row 22: H22 1449.3334->1452.8572, J22 1674.3->1704.7778, L22 1674.3->1704.7778, N22 -2264.3->-2294.7778
$ws.Range("H22").Value = 1452.8572
$ws.Range("J22").Value = 1704.7778
$ws.Range("L22").Value = 1704.7778
$ws.Range("N22").Value = -2294.7778

# Row 27: H27 1449.3334->1452.8572, J27 1674.3->1704.7778, L27 1674.3->1704.7778, N27 -1888.3->-1918.7778
$ws.Range("H27").Value = 1452.8572
$ws.Range("J27").Value = 1704.7778
$ws.Range("L27").Value = 1704.7778
$ws.Range("N27").Value = -1918.7778

# Row 46: H46 1453.3125->1476.9333, J46 1460.9286->1488.7693, L46 1460.9286->1488.7693, N46 -1836.9286->-1864.7693
$ws.Range("H46").Value = 1476.9333
$ws.Range("J46").Value = 1488.7693
$ws.Range("L46").Value = 1488.7693
$ws.Range("N46").Value = -1864.7693

# Row 55: H55 610.44446->609.55554, I55 733.3333->586.25, J55 549->628.2, K55 733.3333->586.25, L55 549->628.2, M55 -560.3333->-413.25, N55 -895->-974.2
$ws.Range("H55").Value = 609.55554
$ws.Range("I55").Value = 586.25
$ws.Range("J55").Value = 628.2
$ws.Range("K55").Value = 586.25
$ws.Range("L55").Value = 628.2
$ws.Range("M55").Value = -413.25
$ws.Range("N55").Value = -974.2

# Row 132: H132 2342.1843->2260.6904, I132 1981.2354->1907.2632, J132 2634.3809->2552.652, K132 5943.706200000001->5721.7896, L132 7903.1427->7657.956, M132 -3413.706200000001->-3191.7896, N132 -12963.1427->-12717.956
$ws.Range("H132").Value = 2260.6904
$ws.Range("I132").Value = 1907.2632
$ws.Range("J132").Value = 2552.652
$ws.Range("K132").Value = 5721.7896
$ws.Range("L132").Value = 7657.956
$ws.Range("M132").Value = -3191.7896
$ws.Range("N132").Value = -12717.956

$ws = $wb.Worksheets.Item("WVR")
# Row 45: H45 74725->36974.75, J45 74725->36974.75, L45 74725->36974.75, N45 -75707->-37956.75
$ws.Range("H45").Value = 36974.75
$ws.Range("J45").Value = 36974.75
$ws.Range("L45").Value = 36974.75
$ws.Range("N45").Value = -37956.75

# Row 62: H62 6049.933->6049.8667, I62 5809.5713->5594.375, J62 6260.25->6570.4287, K62 5809.5713->5594.375, L62 6260.25->6570.4287, M62 -5185.5713->-4970.375, N62 -7508.25->-7818.4287
$ws.Range("H62").Value = 6049.8667
$ws.Range("I62").Value = 5594.375
$ws.Range("J62").Value = 6570.4287
$ws.Range("K62").Value = 5594.375
$ws.Range("L62").Value = 6570.4287
$ws.Range("M62").Value = -4970.375
$ws.Range("N62").Value = -7818.4287

# Row 65: H65 6049.933->6049.8667, I65 5809.5713->5594.375, J65 6260.25->6570.4287, K65 29047.8565->27971.875, L65 31301.25->32852.14350000001, M65 -25927.8565->-24851.875, N65 -37541.25->-39092.14350000001
$ws.Range("H65").Value = 6049.8667
$ws.Range("I65").Value = 5594.375
$ws.Range("J65").Value = 6570.4287
$ws.Range("K65").Value = 27971.875
$ws.Range("L65").Value = 32852.14350000001
$ws.Range("M65").Value = -24851.875
$ws.Range("N65").Value = -39092.14350000001

# Row 81: H81 5527.75->4865.9287, I81 6133.8->5260.6665, K81 12267.6->10521.333, M81 -11206.6->-9460.333000000001
$ws.Range("H81").Value = 4865.9287
$ws.Range("I81").Value = 5260.6665
$ws.Range("K81").Value = 10521.333
$ws.Range("M81").Value = -9460.333000000001

# Row 84: H84 5527.75->4865.9287, I84 6133.8->5260.6665, K84 61338->52606.665, M84 -56034->-47302.665
$ws.Range("H84").Value = 4865.9287
$ws.Range("I84").Value = 5260.6665
$ws.Range("K84").Value = 52606.665
$ws.Range("M84").Value = -47302.665

# Row 132: H132 2573.0642->2097.4358, I132 2551.0212->2432.8723, J132 2606.484->1588.871, K132 7653.0636->7298.6169, L132 7819.451999999999->4766.613, M132 -5123.0636->-4768.6169, N132 -12879.452->-9826.613000000001
$ws.Range("H132").Value = 2097.4358
$ws.Range("I132").Value = 2432.8723
$ws.Range("J132").Value = 1588.871
$ws.Range("K132").Value = 7298.6169
$ws.Range("L132").Value = 4766.613
$ws.Range("M132").Value = -4768.6169
$ws.Range("N132").Value = -9826.613000000001

